$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '68.993.85'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -1.12%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.506.81'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -1.89%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '572.10'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.93%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '184.27'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -2.69%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.615'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -2.73%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '3.500.13'
$cell.Style = "Normal"

$ws.Range("E9").Value = '  +0.07%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.187'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +4.88%  '

$ws.Range("E11").Value = '  -2.12%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '54.16'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -2.92%  '

$ws.Range("E13").Value = '  -0.26%  '

$ws.Range("E14").Value = '  -2.10%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '4.068.60'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.95%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '19.30'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -2.36%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '68.907.56'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.17%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.490.85'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.32%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '12.26'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -3.16%  '

$ws.Range("E20").Value = '  -1.21%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '541.83'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +14.16%  '

$ws.Range("E22").Value = '  -2.83%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '18.78'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -3.00%  '

$ws.Range("E24").Value = '  -0.61%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '4.41'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.44%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '94.04'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.91%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '2.92'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -3.01%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.80'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -2.08%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '9.16'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -2.19%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '31.76'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.97%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '7.25'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -9.26%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '12.56'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.72%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '64.72'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.02%  '

$ws.Range("E34").Value = '  -4.65%  '

$ws.Range("E35").Value = '  -2.79%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '37.90'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.96%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.14%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.397'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '3.01'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +4.52%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0766'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -3.98%  '

$ws.Range("E41").Value = '  -1.33%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.133'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -3.28%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '3.34'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -2.95%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '3.230.21'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.51'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +3.48%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.98'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -3.61%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.0440'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.66%  '

$ws.Range("E48").Value = '  -2.67%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '8.97'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -4.74%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.08%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '137.97'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +2.51%  '
